# Update "想去人数" (want-to-go count) figures in column F across the
# workbook's sheets, matching the refreshed output generated at 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 1744
$ws1.Range("F6").Value  = 3313
$ws1.Range("F7").Value  = 1002
$ws1.Range("F8").Value  = 2170
$ws1.Range("F9").Value  = 2084
$ws1.Range("F11").Value = 591
$ws1.Range("F13").Value = 1654
$ws1.Range("F14").Value = 381
$ws1.Range("F18").Value = 186
$ws1.Range("F19").Value = 1547
$ws1.Range("F20").Value = 596
$ws1.Range("F21").Value = 699
$ws1.Range("F22").Value = 580
$ws1.Range("F23").Value = 12133
$ws1.Range("F24").Value = 12153
$ws1.Range("F25").Value = 896
$ws1.Range("F28").Value = 15
$ws1.Range("F30").Value = 327
$ws1.Range("F31").Value = 1905
$ws1.Range("F33").Value = 551

# --- Sheet "演出" (performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 118
$ws2.Range("F6").Value = 39

# --- Sheet "全部类型" (all types, merged view) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value  = 1744
$ws4.Range("F7").Value  = 3313
$ws4.Range("F8").Value  = 1002
$ws4.Range("F9").Value  = 2170
$ws4.Range("F10").Value = 2084
$ws4.Range("F12").Value = 591
$ws4.Range("F14").Value = 1654
$ws4.Range("F15").Value = 381
$ws4.Range("F22").Value = 186
$ws4.Range("F23").Value = 1547
$ws4.Range("F24").Value = 596
$ws4.Range("F25").Value = 699
$ws4.Range("F26").Value = 580
$ws4.Range("F27").Value = 12133
$ws4.Range("F28").Value = 12153
$ws4.Range("F29").Value = 896
$ws4.Range("F32").Value = 15
$ws4.Range("F34").Value = 327
$ws4.Range("F35").Value = 1905
$ws4.Range("F36").Value = 118
$ws4.Range("F37").Value = 39
$ws4.Range("F39").Value = 551
